# edit.ps1 - applies the lab04 report edits described by the commit
# "feat(main): add files lab-4" to $word.ActiveDocument.
#
# Summary of edits (see task diff):
#  1. Theoretical-intro paragraph "///" -> full sentence incl. a quoted word.
#  2. Fix typo "работе и 4 и" -> "работе 4 и" (drop stray "и").
#  3/4. Figure 2 alt-text + caption: "рис.2" -> real description.
#  5. Replace unresolved cross-reference "[-fig:003]" with "3".
#  6/7. Figure 3 alt-text + caption: "рис.3" -> real description (with 'make').
#  8. Replace unresolved cross-reference "[-fig:004]" with "4".
#  9/10. Figure 4 alt-text + caption: "рис.4" -> real description.
#  11. Replace unresolved cross-reference "[-fig:005]" with a bold
#      "¿fig:005?" placeholder run, keeping the parens as separate runs.
#  12. Figure 5 alt-text: "рис.5" -> real description.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "///" -> "Язык разметки Markdown является "легковесным", при чем
#    данный термин обязательно подлежит заключению в кавычки"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "///", $true, $false, $false, $false, $false, $true, 1, $false,
    [char]0x201C + "легковесным" + [char]0x201D,
    2) | Out-Null

$d.Content.Find.Execute(
    [char]0x201C + "легковесным" + [char]0x201D, $true, $false, $false, $false, $false, $true, 1, $false,
    "Язык разметки Markdown является " + [char]0x201C + "легковесным" + [char]0x201D + ", при чем данный термин обязательно подлежит заключению в кавычки",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2) drop the stray "и" -> "...работе и 4 и..." -> "...работе 4 и..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "лабораторной работе и 4 и провожу", $true, $false, $false, $false, $false, $true, 1, $false,
    "лабораторной работе 4 и провожу",
    2) | Out-Null

# ---------------------------------------------------------------------
# 3/4) Figure 2 - alt text & caption
# ---------------------------------------------------------------------
$d.InlineShapes.Item(2).AlternativeText = "Рис. 2: Компиляция шаблона с использованием Makefile"

$d.Content.Find.Execute(
    "Рис. 2: рис.2", $true, $false, $false, $false, $false, $true, 1, $false,
    "Рис. 2: Компиляция шаблона с использованием Makefile",
    2) | Out-Null

# ---------------------------------------------------------------------
# 5) "[-fig:003]" -> "3"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "(рис. [-fig:003])", $true, $false, $false, $false, $false, $true, 1, $false,
    "(рис. 3)",
    2) | Out-Null

# ---------------------------------------------------------------------
# 6/7) Figure 3 - alt text & caption (caption ends with 'make' in quotes)
# ---------------------------------------------------------------------
$d.InlineShapes.Item(3).AlternativeText = "Рис. 3: Результат работы команды " + [char]0x2018 + "make" + [char]0x2019

$d.Content.Find.Execute(
    "Рис. 3: рис.3", $true, $false, $false, $false, $false, $true, 1, $false,
    "Рис. 3: Результат работы команды " + [char]0x2018 + "make" + [char]0x2019,
    2) | Out-Null

# ---------------------------------------------------------------------
# 8) "[-fig:004]" -> "4"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "(рис. [-fig:004])", $true, $false, $false, $false, $false, $true, 1, $false,
    "(рис. 4)",
    2) | Out-Null

# ---------------------------------------------------------------------
# 9/10) Figure 4 - alt text & caption
# ---------------------------------------------------------------------
$d.InlineShapes.Item(4).AlternativeText = "Рис. 4: Удаление файлов с помощью Makefile"

$d.Content.Find.Execute(
    "Рис. 4: рис.4", $true, $false, $false, $false, $false, $true, 1, $false,
    "Рис. 4: Удаление файлов с помощью Makefile",
    2) | Out-Null

# ---------------------------------------------------------------------
# 11) "(рис. [-fig:005])" -> "(рис. ¿fig:005?)" with the placeholder bold
# ---------------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute(
    "(рис. [-fig:005])", $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0)
if ($found) {
    $find.Text = "(рис. )"
    $insPoint = $find.Start + 6
    $r = $d.Range($insPoint, $insPoint)
    $r.InsertBefore([char]0xBF + "fig:005" + [char]0x3F)
    $boldRange = $d.Range($insPoint, $insPoint + 9)
    $boldRange.Bold = 1
}

# ---------------------------------------------------------------------
# 12) Figure 5 - alt text
# ---------------------------------------------------------------------
$d.InlineShapes.Item(5).AlternativeText = "Проверка отсутствия удаленных файлов"

Write-Output "done"
